$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide the rows that correspond to "Week 1" entries already paid out (rows collapsed in the board).
$hiddenRows = @(5, 7, 9, 10, 12, 15, 17, 20, 22, 25, 27, 30, 33, 35, 38, 40, 41, 43)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}

# Update the active selection / view to D45.
$ws.Range("D45").Select()
